# Apply the "add 2022-Q3 data" edit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet: shift existing quarter rows down
#    and insert the new 2022-Q3 summary row at the top of the data.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Copy the formatting of row 3 (A3, bold/border style) down into row 4
# before writing the new 2022-Q1 row there, so the new row keeps the
# same look as the other data rows.
$totals.Range("A3").Copy()
$totals.Range("A4").PasteSpecial(-4122)

# Row 4 becomes what used to be row 3 (2022-Q1 / 4 / 1.05)
$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2022-Q1"
$totals.Range("C4").Value = 4
$totals.Range("D4").Value = 1.05

# Row 3 becomes what used to be row 2 (2022-Q2 / 4 / 1.3)
$totals.Range("B3").Value = "2022-Q2"
$totals.Range("C3").Value = 4
$totals.Range("D3").Value = 1.3

# Row 2 becomes the new 2022-Q3 summary
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 6
$totals.Range("D2").Value = 0.66

# ---------------------------------------------------------------------
# 2. Insert a brand-new worksheet "2022-Q3" right before the existing
#    "2022-Q2" sheet, so the final sheet order is:
#    总计, 2022-Q3, 2022-Q2, 2022-Q1
# ---------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($existingQ2)
$q3.Name = "2022-Q3"

# Header row text
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Copy the header style (bold + border, matching the style used
# elsewhere in this workbook) from the 总计 sheet's header row, and
# apply it across the full header row.
$totals.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# Copy the column-A "index" style (bold + border) from 总计!A2.
$totals.Range("A2").Copy()
$q3.Range("A2:A7").PasteSpecial(-4122)

$data = @(
    @(0, "001877", "宝盈国家安全沪港深股票A",   "8.74", "91.81", "3.27", "0.2858", 9),
    @(1, "013895", "宝盈成长精选混合A",         "8.51", "90.34", "2.96", "0.2519", 7),
    @(2, "013896", "宝盈成长精选混合C",         "2.93", "90.34", "2.96", "0.0867", 7),
    @(3, "090019", "大成景恒混合A",             "1.13", "93.98", "1.68", "0.0190", 9),
    @(4, "006038", "大成景恒混合C",             "0.45", "93.98", "1.68", "0.0076", 9),
    @(5, "013613", "宝盈国家安全沪港深股票C",   "0.23", "91.81", "3.27", "0.0075", 9)
)

$rowNum = 2
foreach ($rec in $data) {
    # Column A: numeric row index (0-based), bold/bordered style already applied above.
    $q3.Cells.Item($rowNum, 1).Value = $rec[0]

    # Column B: fund code - leading apostrophe forces text storage so
    # leading zeros (e.g. "001877") are preserved instead of being
    # read as the number 1877.
    $q3.Cells.Item($rowNum, 2).Value = "'" + $rec[1]

    # Column C: fund name - plain Chinese text, no numeric conversion risk.
    $q3.Cells.Item($rowNum, 3).Value = $rec[2]

    # Columns D-G: decimal-looking figures stored as text in the source
    # data (matches the style of the other quarter sheets), so again
    # force text via the leading apostrophe.
    $q3.Cells.Item($rowNum, 4).Value = "'" + $rec[3]
    $q3.Cells.Item($rowNum, 5).Value = "'" + $rec[4]
    $q3.Cells.Item($rowNum, 6).Value = "'" + $rec[5]
    $q3.Cells.Item($rowNum, 7).Value = "'" + $rec[6]

    # Column H: numeric rank.
    $q3.Cells.Item($rowNum, 8).Value = $rec[7]

    $rowNum = $rowNum + 1
}
